# Simplified -> Traditional Chinese localisation pass for
# "Removal of USDT Tether Omni - Reminder email to clients.docx" (zh)
#
# Each call does a literal (non-wildcard) Find/Replace over the whole
# document body. Every source string below is unique in the document, so
# there is no risk of an unintended match elsewhere.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute(
        $find,      # FindText
        $true,      # MatchCase
        $false,     # MatchWholeWord
        $false,     # MatchWildcards
        $false,     # MatchSoundsLike
        $false,     # MatchAllWordForms
        $true,      # Forward
        1,          # Wrap (wdFindContinue)
        $false,     # Format
        $replace,   # ReplaceWith
        2           # Replace (wdReplaceAll)
    ) | Out-Null
}

# Title / header line
Replace-Text "向 ROW 客户发送提醒电子邮件" "向 ROW 客戶傳送提醒電子郵件"

# Subject line
Replace-Text "将于 9 月 29 日移除 Tether Omni (USDT)" "將於 9 月 29 日移除 Tether Omni (USDT)"

# Body heading
Replace-Text "向 Tether Omni 说再见" "向 Tether Omni 道別"

# Intro paragraph
Replace-Text "自 2023 年 9 月 29 日格林威治标准时间 00:00 起，Deriv 将停止提供Tether Omni (USDT) 作为账户货币。 这是因为 Tether 已停止支持 USDT 的 Omni 转账。" "自 2023 年 9 月 29 日格林威治標準時間 00:00 起，Deriv 將停止提供Tether Omni (USDT) 作為帳戶貨幣。 這是因為 Tether 已停止支援 USDT 的 Omni 轉帳。"

# "What do you need to do?" heading
Replace-Text "需要做什么？" "需要做什麼？"

# "If the USDT account "
Replace-Text "如果 USDT 账户 " "如果 USDT 帳戶 "

# remainder after [account ID]
Replace-Text " 中有余额，请在上述日期之前提取余额。 如果有持仓头寸，提取余额之前请先平仓。" " 中有餘額，請在上述日期之前提取餘額。 若有持倉頭寸，提取餘額前請先平倉。"

# Button text
Replace-Text "查看账户" "檢查帳戶"

# "Important notice" paragraph
Replace-Text "USDT 账户将于 2023 年 9 月 29 日格林尼治标准时间 00:00 关闭。 任何持仓头寸将在上述日期后自动平仓，账户余额将转移到最后活跃的账户" "USDT 帳戶將於 2023 年 9 月 29 日格林尼治標準時間 00:00 關閉。 任何持倉頭寸將在上述日期後自動平倉，帳戶餘額將轉移到最後活躍的帳戶"

# Trailing sentence after the comment markers
Replace-Text "在此过程中将采用标准汇率和费用。" "在此過程中將採用標準匯率和費用。"

# "If you have any questions" line
Replace-Text "如有任何疑问，请通过以下方式联系我们：" "如有任何疑問，請透過以下方式聯繫我們："

# "Live chat" link text
Replace-Text "实时聊天" "即時聊天"
